# The fcs processing equation moved from plate_01 to plate_02 (see commit
# message: "fcs file now in plate 2 of example.xlsx"), and the
# "Transformations" sheet became the active tab/selection in the workbook
# (previously "Samples" was the active tab).

$wb = $excel.ActiveWorkbook

# Update the transformation formula to reference plate_02 instead of plate_01.
$transformations = $wb.Worksheets.Item("Transformations")
$transformations.Range("B2").Value = 'process_fcs("plate_02",["FSC","SSC"],["FL1"])'

# Move the active sheet / selection to Transformations (was Samples before).
$transformations.Activate() | Out-Null
$transformations.Range("B3").Select() | Out-Null

Write-Output "Updated Transformations!B2 and switched active sheet to Transformations"
